$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.855.92'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.636.22'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.51'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('E6').Value = '  -0.42%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.15'
$ws.Range('E10').Value = '  +4.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.667.47'
$ws.Range('E12').Value = '  +2.60%  '
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.862.07'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.26'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.864.70'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.94'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('E21').Value = '  +0.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.93'
$ws.Range('E22').Value = '  +1.46%  '
$ws.Range('E23').Value = '  +3.35%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  -3.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '138.47'
$ws.Range('E26').Value = '  -1.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.122'
$ws.Range('E27').Value = '  -4.79%  '
$ws.Range('E28').Value = '  +1.50%  '
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('E31').Value = '  +1.17%  '
$ws.Range('E32').Value = '  +0.59%  '
$ws.Range('E33').Value = '  +1.60%  '
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('E35').Value = '  +0.64%  '
$ws.Range('E36').Value = '  +1.02%  '
$ws.Range('E37').Value = '  +1.55%  '
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.121.84'
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0158'
$ws.Range('E40').Value = '  +1.66%  '
$ws.Range('E41').Value = '  +0.57%  '
$ws.Range('E42').Value = '  -1.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.47'
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('E44').Value = '  +0.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₆0112'
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.47'
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('E47').Value = '  -4.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0504'
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.66'
$ws.Range('E49').Value = '  +1.35%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('E51').Value = '  +0.06%  '
